$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -13.01
$ws.Range("C4").Value  = -13.318
$ws.Range("D6").Value  = -7.882000000000001
$ws.Range("C7").Value  = -13.138
$ws.Range("D7").Value  = -7.334999999999999
$ws.Range("C8").Value  = -12.672
$ws.Range("D8").Value  = -7.858
$ws.Range("B11").Value = 6.025
$ws.Range("B12").Value = 5.568
$ws.Range("C12").Value = -13.213
$ws.Range("C14").Value = -12.081
$ws.Range("B15").Value = 6.845000000000001
$ws.Range("D19").Value = -7.790999999999999
$ws.Range("D21").Value = -7.858
$ws.Range("C22").Value = -13.005
$ws.Range("D24").Value = -7.934
$ws.Range("D25").Value = -7.861
